# Update the "address" column (D) values on the "Child" worksheet to the
# newly generated child-travel coordinates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$values = @{
    2  = "-4.18,-8.88"
    3  = "4.52,-9.26"
    4  = "-1.98,-2.1"
    5  = "-1.65,-8.14"
    6  = "-0.28,2.19"
    7  = "-6.55,4.12"
    8  = "-9.79,-8.09"
    9  = "-7.11,9.53"
    10 = "5.4,-6.02"
    11 = "-8.77,7.51"
    12 = "8.81,2.38"
    13 = "9.63,4.02"
    14 = "1.45,4.78"
    15 = "0.33,2.02"
    16 = "3.37,7.62"
    17 = "9.71,4.53"
    18 = "6.19,7.48"
    19 = "-3.64,7.49"
    20 = "9.32,9.44"
    21 = "6.54,0.52"
    22 = "9.6,-1.85"
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
